# ---------------------------------------------------------------------------
# Adds 20 new LDAP/User entries to the "MAIN" sheet (rows 595-614), mirrors
# them into the computed "DataLookUp" sheet, moves the trailing blank
# template row from 595 -> 615 on both sheets, widens column C (MAIN) /
# column D (DataLookUp) to fit the longer LDAP strings, extends the
# conditional-formatting range + duplicates its dxf, and switches the
# active sheet/tab back to MAIN.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$main = $wb.Worksheets.Item("MAIN")
$lookup = $wb.Worksheets.Item("DataLookUp")

# New (LDAP) / User_ pairs, in the order they were appended.
$pairs = @(
    @("(LDAP) ridwan.nur", "User_RidwanNurMutaqqin"),
    @("(LDAP) muhammad.rizal", "User_MuhammadRizal"),
    @("(LDAP) andi.mapasoro", "User_AndiMapasoro"),
    @("(LDAP) taufik.aminudin", "User_TaufikAminudin"),
    @("(LDAP) taufiq.aminudin", "User_TaufiqAminudin"),
    @("(LDAP) risty.alawiyah", "User_RistyAlawiyah"),
    @("(LDAP) agus.sobar", "User_AgusSobar"),
    @("(LDAP) reza.sanjaya", "User_RezaSanjaya"),
    @("(LDAP) kiki.prasetyo", "User_KikiPrasetyo"),
    @("(LDAP) ulul.azmi", "User_MochamadUlulAzmiWiraUtama"),
    @("(LDAP) arie.gustaman", "User_ArieGustaman"),
    @("(LDAP) zulfikar.siregar", "User_ZulfikarKSiregar"),
    @("(LDAP) deby.syahputra", "User_DebySyahputa"),
    @("(LDAP) Risaldi", "User_Risaldi"),
    @("(LDAP) muhammad.hairul", "User_MuhammadHairul"),
    @("(LDAP) ibrahim.arya", "User_IbrahimAryaYuda"),
    @("(LDAP) m.bilal", "User_MuhammadBilalSharony"),
    @("(LDAP) jamil", "User_Jamil"),
    @("(LDAP) wiyadianto", "User_Wiyadianto"),
    @("(LDAP) rahmad.dwifar", "User_RahmadDwifar")
)

$n = $pairs.Length
$oldLast = 594        # last populated data row before the edit
$oldBlank = 595        # old trailing blank template row
$newBlank = $oldBlank + $n   # 615 - trailing blank template row after the edit

# ---------------------------------------------------------------------------
# 1) MAIN sheet - move the blank template row's formatting out of the way
#    first (copy row 594's "real data" look down onto row 615), then
#    restyle row 594 itself like the old blank row, matching what Excel
#    does when a block is inserted just above the sheet's final row.
# ---------------------------------------------------------------------------

# Stash the two "looks" we need (the real-data look from row 594 and the
# blank-row look from row 595) off in scratch cells first, since the
# paste operations below overwrite both source rows before all the
# copies that need them are done.
$main.Range("B" + $oldLast + ":D" + $oldLast).Copy()
$main.Range("B1000:D1000").PasteSpecial(-4122)
$main.Range("F" + $oldLast + ":G" + $oldLast).Copy()
$main.Range("F1000:G1000").PasteSpecial(-4122)

$main.Range("B" + $oldBlank + ":D" + $oldBlank).Copy()
$main.Range("B1001:D1001").PasteSpecial(-4122)
# NB: F on the old blank row (595) has no cell of its own (it's blank,
# implicitly styled by the column default) - use F594's *numeric* look
# for the new blank row's F cell instead, matching the target (F615 is
# styled/formulaed the same as the other F-column cells).
$main.Range("F" + $oldLast).Copy()
$main.Range("F1001").PasteSpecial(-4122)

# New data rows 595..614 get the same look the old last data row (594) had.
$main.Range("B1000:D1000").Copy()
$main.Range("B595:D" + ($oldBlank + $n - 1)).PasteSpecial(-4122)
$main.Range("F1000:G1000").Copy()
$main.Range("F595:G" + ($oldBlank + $n - 1)).PasteSpecial(-4122)

# The trailing blank row moves from 595 to 615, carrying its look with it.
$main.Range("B1001:D1001").Copy()
$main.Range("B" + $newBlank + ":D" + $newBlank).PasteSpecial(-4122)
$main.Range("F1001").Copy()
$main.Range("F" + $newBlank).PasteSpecial(-4122)

# Row 594 (no longer the last data row) now looks like a blank row.
$main.Range("B1001:D1001").Copy()
$main.Range("B" + $oldLast + ":D" + $oldLast).PasteSpecial(-4122)

# Clean up the scratch cells.
$main.Range("B1000:G1001").ClearFormats()
$main.Range("B1000:G1001").ClearContents()

$main.Cells.Item(1, 1).Select() | Out-Null

for ($i = 0; $i -lt $n; $i++) {
    $r = 595 + $i
    $ldap = $pairs[$i][0]
    $user = $pairs[$i][1]

    $main.Cells.Item($r, 2).Value = $ldap
    $main.Cells.Item($r, 3).Value = $user
    $main.Cells.Item($r, 6).Formula = '=F' + ($r - 1) + ' + IF(EXACT(G' + $r + ', ""), 0, 1)'
    $main.Cells.Item($r, 7).Formula = '=IF(EXACT(B' + $r + ', ""), "", CONCATENATE("PERFORM ""SchSysConfig"".""Func_TblDBObject_User_SET""(varSystemLoginSession, varInstitutionBranchID, varBaseCurrencyID, ''", B' + $r + ', "'', ", IF(EXACT(C' + $r + ', ""), "null", CONCATENATE("''", C' + $r + ', "''")), ", ", IF(EXACT(D' + $r + ', ""), "''''", CONCATENATE("''", D' + $r + ', "''")), ");"))'
}

# Trailing blank row (615): clear any stray values, keep only the rolling
# F counter formula (same shape as the old blank row 595).
$main.Cells.Item($newBlank, 2).ClearContents()
$main.Cells.Item($newBlank, 3).ClearContents()
$main.Cells.Item($newBlank, 7).ClearContents()
$main.Cells.Item($newBlank, 6).Formula = '=F' + ($newBlank - 1) + ' + IF(EXACT(G' + $newBlank + ', ""), 0, 1)'

# ---------------------------------------------------------------------------
# 2) DataLookUp sheet - same row shuffle, but the cells are themselves
#    formulas mirroring MAIN, so there is no literal value to type in.
# ---------------------------------------------------------------------------

$lookup.Range("B" + $oldLast + ":D" + $oldLast).Copy()
$lookup.Range("B1000:D1000").PasteSpecial(-4122)
$lookup.Range("B" + $oldBlank + ":D" + $oldBlank).Copy()
$lookup.Range("B1001:D1001").PasteSpecial(-4122)

$lookup.Range("B1000:D1000").Copy()
$lookup.Range("B595:D" + ($oldBlank + $n - 1)).PasteSpecial(-4122)

$lookup.Range("B1001:D1001").Copy()
$lookup.Range("B" + $newBlank + ":D" + $newBlank).PasteSpecial(-4122)

$lookup.Range("B1001:D1001").Copy()
$lookup.Range("B" + $oldLast + ":D" + $oldLast).PasteSpecial(-4122)

$lookup.Range("B1000:D1001").ClearFormats()
$lookup.Range("B1000:D1001").ClearContents()

for ($i = 0; $i -lt $n; $i++) {
    $r = 595 + $i
    $lookup.Cells.Item($r, 2).Formula = '=IF(EXACT(MAIN!$G' + $r + ', ""), "", MAIN!$F' + $r + ')'
    $lookup.Cells.Item($r, 3).Formula = '=IF(EXACT(MAIN!$G' + $r + ', ""), "", MAIN!$B' + $r + ')'
    $lookup.Cells.Item($r, 4).Formula = '=IF(EXACT(MAIN!$G' + $r + ', ""), "", MAIN!$C' + $r + ')'
}

$lookup.Cells.Item($newBlank, 2).Formula = '=IF(EXACT(MAIN!$G' + $newBlank + ', ""), "", MAIN!$F' + $newBlank + ')'
$lookup.Cells.Item($newBlank, 3).Formula = '=IF(EXACT(MAIN!$G' + $newBlank + ', ""), "", MAIN!$B' + $newBlank + ')'
$lookup.Cells.Item($newBlank, 4).Formula = '=IF(EXACT(MAIN!$G' + $newBlank + ', ""), "", MAIN!$C' + $newBlank + ')'

# ---------------------------------------------------------------------------
# 3) Column widths - column C (MAIN) / column D (DataLookUp) now need to fit
#    the longest new LDAP string.
# ---------------------------------------------------------------------------

$main.Columns.Item(3).ColumnWidth = 26.140625
$lookup.Columns.Item(4).ColumnWidth = 26.140625

# ---------------------------------------------------------------------------
# 4) Conditional formatting on MAIN!F4:F615 - duplicate the existing dxf
#    (a new FormatCondition is added referencing a *new* dxf slot) and
#    drop the old one so only the newly-added rule covers the sheet.
# ---------------------------------------------------------------------------

$cf = $main.Range("F4:F594").FormatConditions.Item(1)
$main.Range("F4:F615").FormatConditions.Add(2, 3, "=EXACT(F3,F4)") | Out-Null
$newCf = $main.Range("F4:F615").FormatConditions.Item($main.Range("F4:F615").FormatConditions.Count)
$newCf.Interior.ThemeColor = 1
$cf.Delete()

# ---------------------------------------------------------------------------
# 5) Switch the active sheet back to MAIN (was DataLookUp before the edit).
# ---------------------------------------------------------------------------

$main.Activate()
$main.Range("F619").Select() | Out-Null
$lookup.Range("C608").Select() | Out-Null
$main.Activate()
